$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New column I ("hsphere_plan_id"): header + values for the existing rows.
# Copy the format from column H so the new cells inherit the same style
# already used across row 1-3, then overwrite with the new content.
# ---------------------------------------------------------------------------
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "hsphere_plan_id"

$ws.Range("H2").Copy() | Out-Null
$ws.Range("I2").PasteSpecial(-4122) | Out-Null
$ws.Range("I2").Value = 1234.0

$ws.Range("H3").Copy() | Out-Null
$ws.Range("I3").PasteSpecial(-4122) | Out-Null
$ws.Range("I3").Value = 2112.0

$ws.Range("I1").ColumnWidth = 15.25

# ---------------------------------------------------------------------------
# Row 2: rename the plan and bump its opspi_account_id.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "hsphere_user_plan_2"
$ws.Range("H2").Value = 5.0

# ---------------------------------------------------------------------------
# Row 3: rename the plan, fix plan_duration, bump opspi_account_id.
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "hsphere_user_plan_3"
$ws.Range("E3").Value = 1.0
$ws.Range("H3").Value = 5.0

# ---------------------------------------------------------------------------
# Append rows 4-7 (row 4 replaces the old hsphere_plan_12 data). Formats are
# copied down from row 3 / column I first, so the new rows start out on the
# same style as the rest of the table before the new font is applied below.
# ---------------------------------------------------------------------------
$ws.Range("A3:H3").Copy() | Out-Null
$ws.Range("A4:H7").PasteSpecial(-4122) | Out-Null

$ws.Range("I3").Copy() | Out-Null
$ws.Range("I4:I7").PasteSpecial(-4122) | Out-Null

$ws.Range("A4").Value = "hsphere_user_plan_4"
$ws.Range("B4").Value = 100.0
$ws.Range("C4").Value = 500.0
$ws.Range("D4").Value = 1000.0
$ws.Range("E4").Value = 1.0
$ws.Range("F4").Value = $true
$ws.Range("G4").Value = $true
$ws.Range("H4").Value = 5.0
$ws.Range("I4").Value = 3542.0

$ws.Range("A5").Value = "hsphere_user_plan_5"
$ws.Range("B5").Value = 50.0
$ws.Range("C5").Value = 300.0
$ws.Range("D5").Value = 600.0
$ws.Range("E5").Value = 1.0
$ws.Range("F5").Value = $true
$ws.Range("G5").Value = $true
$ws.Range("H5").Value = 5.0
$ws.Range("I5").Value = 4532.0

$ws.Range("A6").Value = "hsphere_user_plan_6"
$ws.Range("B6").Value = 100.0
$ws.Range("C6").Value = 500.0
$ws.Range("D6").Value = 1000.0
$ws.Range("E6").Value = 1.0
$ws.Range("F6").Value = $true
$ws.Range("G6").Value = $true
$ws.Range("H6").Value = 5.0
$ws.Range("I6").Value = 5643.0

$ws.Range("A7").Value = "hsphere_user_plan_7"
$ws.Range("B7").Value = 50.0
$ws.Range("C7").Value = 300.0
$ws.Range("D7").Value = 600.0
$ws.Range("E7").Value = 1.0
$ws.Range("F7").Value = $true
$ws.Range("G7").Value = $true
$ws.Range("H7").Value = 5.0
$ws.Range("I7").Value = 6787.0

# ---------------------------------------------------------------------------
# New rows get their own font (Arial / theme text color, no inherited minor
# scheme) plus per-column alignment: numerics right-aligned, booleans
# centered, plan names left as-is.
# ---------------------------------------------------------------------------
$ws.Range("A4:I7").Font.Name = "Arial"
$ws.Range("B4:E7").HorizontalAlignment = -4152
$ws.Range("H4:H7").HorizontalAlignment = -4152
$ws.Range("I4:I7").HorizontalAlignment = -4152
$ws.Range("F4:G7").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Rows 6-7 also carry a run of empty, styled filler cells out to column AB
# (matches the wide "blank trailer" written for these two rows).
# ---------------------------------------------------------------------------
$ws.Range("J6:AB7").Value = ""
$ws.Range("J6:AB7").Font.Name = "Arial"
$ws.Range("J6:AB7").Font.ThemeColor = 1
